$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.594.21'
$ws.Range('E2').Value = '  +0.47%  '

$ws.Range('D3').Value = '1.920.76'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').Value = '247.28'
$ws.Range('E5').Value = '  +3.01%  '

$ws.Range('E6').Value = '  -0.23%  '

$ws.Range('D7').Value = '0.4744'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').Value = '0.2889'
$ws.Range('E8').Value = '  +1.50%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06830'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.95%  '

$ws.Range('D10').Value = '105.02'
$ws.Range('E10').Value = '  -0.45%  '

$ws.Range('E11').Value = '  -4.59%  '

$ws.Range('D12').Value = '1.918.97'
$ws.Range('E12').Value = '  -0.22%  '

$ws.Range('D13').Value = '0.07703'
$ws.Range('E13').Value = '  +1.53%  '

$ws.Range('D14').Value = '5.284'
$ws.Range('E14').Value = '  +3.07%  '

$ws.Range('D15').Value = '0.6689'
$ws.Range('E15').Value = '  +2.70%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '291.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.70%  '

$ws.Range('D17').Value = '30.595.00'
$ws.Range('E17').Value = '  +0.44%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007600'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.39%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.17%  '

$ws.Range('D20').Value = '12.94'
$ws.Range('E20').Value = '  -0.25%  '

$ws.Range('D21').Value = '5.541'
$ws.Range('E21').Value = '  +6.52%  '

$ws.Range('D22').Value = '2.167.08'
$ws.Range('E22').Value = '  -0.64%  '

$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').Value = '6.376'
$ws.Range('E24').Value = '  +1.14%  '

$ws.Range('D25').Value = '9.408'
$ws.Range('E25').Value = '  +1.23%  '

$ws.Range('E26').Value = '  +1.54%  '

$ws.Range('D27').Value = '21.09'
$ws.Range('E27').Value = '  +7.32%  '

$ws.Range('D28').Value = '2.113'
$ws.Range('E28').Value = '  +3.87%  '

$ws.Range('D29').Value = '0.1068'
$ws.Range('E29').Value = '  -5.09%  '

$ws.Range('D30').Value = '1.397'
$ws.Range('E30').Value = '  +3.58%  '

$ws.Range('D31').Value = '4.176'
$ws.Range('E31').Value = '  +1.48%  '

$ws.Range('D32').Value = '4.068'
$ws.Range('E32').Value = '  +3.64%  '

$ws.Range('D33').Value = '0.05025'
$ws.Range('E33').Value = '  +0.09%  '

$ws.Range('D34').Value = '0.7378'
$ws.Range('E34').Value = '  +0.06%  '

$ws.Range('E35').Value = '  -0.25%  '

$ws.Range('D36').Value = '0.02073'
$ws.Range('E36').Value = '  +5.74%  '

$ws.Range('E37').Value = '  +0.81%  '

$ws.Range('D38').Value = '2.686'

$ws.Range('D39').Value = '2.051'
$ws.Range('E39').Value = '  +1.67%  '

$ws.Range('D40').Value = '111.07'
$ws.Range('E40').Value = '  +3.64%  '

$ws.Range('D41').Value = '0.8777'
$ws.Range('E41').Value = '  +0.45%  '

$ws.Range('D42').Value = '0.4381'
$ws.Range('E42').Value = '  +5.96%  '

$ws.Range('D43').Value = '5.868'
$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('E44').Value = '  -0.20%  '

$ws.Range('D45').Value = '67.54'
$ws.Range('E45').Value = '  -2.07%  '

$ws.Range('D46').Value = '7.245'
$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('D47').Value = '9.284'
$ws.Range('E47').Value = '  +0.81%  '

$ws.Range('D48').Value = '48.38'
$ws.Range('E48').Value = '  +14.65%  '

$ws.Range('E49').Value = '  +1.77%  '

$ws.Range('D50').Value = '34.81'
$ws.Range('E50').Value = '  +0.52%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2480'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.80%  '
